$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 18552.334
$ws.Range("J3").Value = 18552.334
$ws.Range("L3").Value = 18552.334
$ws.Range("N3").Value = -18780.334

$ws.Range("H40").Value = 2120
$ws.Range("I40").Value = 1666.6666
$ws.Range("J40").Value = 2263.158
$ws.Range("K40").Value = 1666.6666
$ws.Range("L40").Value = 2263.158
$ws.Range("M40").Value = -1491.6666
$ws.Range("N40").Value = -2613.158

$ws.Range("H86").Value = 45457084
$ws.Range("I86").Value = 1743
$ws.Range("J86").Value = 200005250
$ws.Range("K86").Value = 1743
$ws.Range("L86").Value = 200005250
$ws.Range("M86").Value = -620
$ws.Range("N86").Value = -200007496

$ws.Range("H89").Value = 45457084
$ws.Range("I89").Value = 1743
$ws.Range("J89").Value = 200005250
$ws.Range("K89").Value = 8715
$ws.Range("L89").Value = 1000026250
$ws.Range("M89").Value = -3099
$ws.Range("N89").Value = -1000037482

$ws.Range("H102").Value = 18552.334
$ws.Range("J102").Value = 18552.334
$ws.Range("L102").Value = 18552.334
$ws.Range("N102").Value = -25042.334

$ws.Range("I106").Value = 4250
$ws.Range("J106").Value = 4500
$ws.Range("K106").Value = 4250
$ws.Range("L106").Value = 4500
$ws.Range("M106").Value = -3619
$ws.Range("N106").Value = -5762

$ws.Range("H132").Value = 1705.3829
$ws.Range("I132").Value = 1726.025
$ws.Range("J132").Value = 1587.4286
$ws.Range("K132").Value = 5178.075000000001
$ws.Range("L132").Value = 4762.2858
$ws.Range("M132").Value = -2648.075000000001
$ws.Range("N132").Value = -9822.2858

$ws.Range("H138").Value = 1785.6
$ws.Range("J138").Value = 1909.9136
$ws.Range("L138").Value = 5729.7408
$ws.Range("N138").Value = -16009.7408

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 47412.316
$ws.Range("I2").Value = 113299.336
$ws.Range("J2").Value = 1798.2307
$ws.Range("K2").Value = 113299.336
$ws.Range("L2").Value = 1798.2307
$ws.Range("M2").Value = -113186.336
$ws.Range("N2").Value = -2024.2307

$ws.Range("H45").Value = 1197.3125
$ws.Range("I45").Value = 1192.909
$ws.Range("J45").Value = 1207
$ws.Range("K45").Value = 1192.909
$ws.Range("L45").Value = 1207
$ws.Range("M45").Value = -815.9090000000001
$ws.Range("N45").Value = -1961

$ws.Range("H102").Value = 3925.7144
$ws.Range("I102").Value = 4296
$ws.Range("K102").Value = 4296
$ws.Range("M102").Value = -2674

$ws.Range("H116").Value = 47412.316
$ws.Range("I116").Value = 113299.336
$ws.Range("J116").Value = 1798.2307
$ws.Range("K116").Value = 113299.336
$ws.Range("L116").Value = 1798.2307
$ws.Range("M116").Value = -111005.336
$ws.Range("N116").Value = -6386.2307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 47412.316
$ws.Range("I3").Value = 113299.336
$ws.Range("J3").Value = 1798.2307
$ws.Range("K3").Value = 113299.336
$ws.Range("L3").Value = 1798.2307
$ws.Range("M3").Value = -113185.336
$ws.Range("N3").Value = -2026.2307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 299.69232
$ws.Range("I7").Value = 345.77777
$ws.Range("J7").Value = 196
$ws.Range("K7").Value = 345.77777
$ws.Range("L7").Value = 196
$ws.Range("M7").Value = -232.77777
$ws.Range("N7").Value = -422

$ws.Range("H31").Value = 2166.4775
$ws.Range("I31").Value = 2281.1458
$ws.Range("J31").Value = 1876.7894
$ws.Range("K31").Value = 2281.1458
$ws.Range("L31").Value = 1876.7894
$ws.Range("M31").Value = -1986.1458
$ws.Range("N31").Value = -2466.7894

$ws.Range("H34").Value = 2166.4775
$ws.Range("I34").Value = 2281.1458
$ws.Range("J34").Value = 1876.7894
$ws.Range("K34").Value = 2281.1458
$ws.Range("L34").Value = 1876.7894
$ws.Range("M34").Value = -2079.1458
$ws.Range("N34").Value = -2280.7894

$ws.Range("H107").Value = 267.5
$ws.Range("I107").Value = 132.5
$ws.Range("J107").Value = 363.92856
$ws.Range("K107").Value = 132.5
$ws.Range("L107").Value = 363.92856
$ws.Range("M107").Value = 1787.5
$ws.Range("N107").Value = -4203.92856

$ws.Range("H132").Value = 1475.9656
$ws.Range("I132").Value = 1012.7917
$ws.Range("J132").Value = 3699.2
$ws.Range("K132").Value = 3038.3751
$ws.Range("L132").Value = 11097.6
$ws.Range("M132").Value = -508.3751000000002
$ws.Range("N132").Value = -16157.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 4769.567
$ws.Range("J96").Value = 4769.567
$ws.Range("L96").Value = 14308.701
$ws.Range("N96").Value = -18426.701

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5676.364
$ws.Range("I113").Value = 7654.2856
$ws.Range("J113").Value = 2215
$ws.Range("K113").Value = 7654.2856
$ws.Range("L113").Value = 2215
$ws.Range("M113").Value = -5484.2856
$ws.Range("N113").Value = -6555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1899.8572
$ws.Range("I7").Value = 2155.3333
$ws.Range("J7").Value = 1440
$ws.Range("K7").Value = 2155.3333
$ws.Range("L7").Value = 1440
$ws.Range("M7").Value = -2043.3333
$ws.Range("N7").Value = -1664

$ws.Range("H16").Value = 5012.269
$ws.Range("I16").Value = 7703.75
$ws.Range("J16").Value = 705.9
$ws.Range("K16").Value = 7703.75
$ws.Range("L16").Value = 705.9
$ws.Range("M16").Value = -7533.75
$ws.Range("N16").Value = -1045.9

$ws.Range("H22").Value = 1465.9445
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1659.1333
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1659.1333
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -2249.1333

$ws.Range("H27").Value = 1465.9445
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1659.1333
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 1659.1333
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -1873.1333

$ws.Range("H40").Value = 632910.1
$ws.Range("I40").Value = 919587.5600000001
$ws.Range("J40").Value = 2219.8
$ws.Range("K40").Value = 919587.5600000001
$ws.Range("L40").Value = 2219.8
$ws.Range("M40").Value = -919451.5600000001
$ws.Range("N40").Value = -2491.8

$ws.Range("H46").Value = 2500.2
$ws.Range("I46").Value = 2500.2
$ws.Range("K46").Value = 2500.2
$ws.Range("M46").Value = -2312.2

$ws.Range("H82").Value = 1772.4736
$ws.Range("I82").Value = 1944
$ws.Range("J82").Value = 1400.8334
$ws.Range("K82").Value = 1944
$ws.Range("L82").Value = 1400.8334
$ws.Range("M82").Value = -1583
$ws.Range("N82").Value = -2122.8334

$ws.Range("H85").Value = 1772.4736
$ws.Range("I85").Value = 1944
$ws.Range("J85").Value = 1400.8334
$ws.Range("K85").Value = 1944
$ws.Range("L85").Value = 1400.8334
$ws.Range("M85").Value = -696
$ws.Range("N85").Value = -3896.8334

$ws.Range("H126").Value = 1899.8572
$ws.Range("I126").Value = 2155.3333
$ws.Range("J126").Value = 1440
$ws.Range("K126").Value = 6465.999899999999
$ws.Range("L126").Value = 4320
$ws.Range("M126").Value = -3995.999899999999
$ws.Range("N126").Value = -9260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 5555
$ws.Range("J97").Value = 5555
$ws.Range("L97").Value = 5555
$ws.Range("N97").Value = -7537

$ws.Range("H113").Value = 491.9091
$ws.Range("I113").Value = 414.13333
$ws.Range("J113").Value = 658.5714
$ws.Range("K113").Value = 1242.39999
$ws.Range("L113").Value = 1975.7142
$ws.Range("M113").Value = 927.6000100000001
$ws.Range("N113").Value = -6315.7142
